# Generate Report for Handoff
# Updates the localization-status workbook with a newly generated handoff
# package (new GUID-named source/target files and refreshed timestamps).

$wb = $excel.ActiveWorkbook

$oldGuid = "b91bc70b-eaea-4e40-84a7-d18b97bf6e8f"
$newGuid = "7eacbddf-b45d-4045-a9d6-3e4957d931f6"

$oldZhXlf = "$oldGuid.14ffc8e0a7885537f7f7aab95f38ddcbef9919ba.zh-cn.xlf"
$newZhXlf = "$newGuid.45c32c21ca1dbe28c7608e8c9b466d773d380610.zh-cn.xlf"

$oldDeXlf = "$oldGuid.14ffc8e0a7885537f7f7aab95f38ddcbef9919ba.de-de.xlf"
$newDeXlf = "$newGuid.45c32c21ca1dbe28c7608e8c9b466d773d380610.de-de.xlf"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/119ea557b937760932f0d78abac9b5d4d899214c/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-02 23:04:33"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-09-02 23:04:29"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-09-02 23:04:33"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")
